# Update cryptocurrency price ("D") and volume-change ("E") columns for rows 2-51.
# Values that look like plain decimals (e.g. "1.003") would otherwise be
# auto-converted to numbers by Excel on assignment (losing trailing zeros / the
# original text formatting), so for those cells we briefly force a Text number
# format before assigning the value, then restore the cell to the Normal style
# so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.529.75"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.912.48"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5260"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.528"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "1.911.58"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06654"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.33%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").Value = "28.628.08"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.697"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.34%  "
$ws.Range("D27").Value = "2.133.65"
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.748"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06736"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02432"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.258"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2225"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.098"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6443"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.189"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6083"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.767"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.37%  "
